$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet holds an "ID" + A/B/C/D/F dataset (rows 2-35).
# Two data rows need to be removed entirely: "RM 232" (row 26) and "SC 92"
# (originally row 28). Deleting row 26 first shifts "SC 92" up to row 27,
# so it is deleted next from that position.
$ws.Rows(26).Delete()
$ws.Rows(27).Delete()

# After the two deletions the remaining rows have shifted up so the
# dimension is now A1:F33. The F column (error/imputed value) also needs
# a handful of per-row corrections on the rows that shifted into their
# new positions (26-33):
#   row 27 -> "SC 101" : F gains a value (17)
#   row 28 -> "SC 105" : F is cleared (now blank)
#   row 29 -> "SC 119" : F is cleared (now blank)
#   row 30 -> "SC 120" : F gains a value (16.89)
#   row 32 -> "SC 193" : F is cleared (now blank)
$ws.Range("F27").Value = 17
$ws.Range("F28").ClearContents()
$ws.Range("F29").ClearContents()
$ws.Range("F30").Value = 16.89
$ws.Range("F32").ClearContents()

Write-Output "edit complete"
